# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" on every sheet
#    that reports a per-locale translation status.
# 2. Shrink the "Status" column(s) that held the old, longer status text
#    so they fit the new, shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
